$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- O4: year header 2021, same style as N4 ---
$ws.Range("N4").Copy()
$ws.Range("O4").PasteSpecial(-4122)
$ws.Range("O4").Value = 2021

# --- O5: Total row, bold numeric style (N5 style + 0.0 number format) ---
$ws.Range("N5").Copy()
$ws.Range("O5").PasteSpecial(-4122)
$ws.Range("O5").Value = 689
$ws.Range("O5").NumberFormat = "0.0"

# --- O9 / O15: dash rows, style based on N9 ("-") + 0.0 number format ---
$ws.Range("N9").Copy()
$ws.Range("O9").PasteSpecial(-4122)
$ws.Range("O9").Value = "-"
$ws.Range("O9").NumberFormat = "0.0"

$ws.Range("N9").Copy()
$ws.Range("O15").PasteSpecial(-4122)
$ws.Range("O15").Value = "-"
$ws.Range("O15").NumberFormat = "0.0"

# --- O16: bottom border row, style based on N16 + 0.0 number format ---
$ws.Range("N16").Copy()
$ws.Range("O16").PasteSpecial(-4122)
$ws.Range("O16").Value = 10.9
$ws.Range("O16").NumberFormat = "0.0"

# --- O6,O7,O8,O10,O11,O12,O13,O14: plain numeric data rows, same style as N10 (already 0.0 fmt) ---
$ws.Range("N10").Copy()

$ws.Range("O6").PasteSpecial(-4122)
$ws.Range("O6").Value = 94.1

$ws.Range("O7").PasteSpecial(-4122)
$ws.Range("O7").Value = 147.1

$ws.Range("O8").PasteSpecial(-4122)
$ws.Range("O8").Value = 10.1

$ws.Range("O10").PasteSpecial(-4122)
$ws.Range("O10").Value = 82.1

$ws.Range("O11").PasteSpecial(-4122)
$ws.Range("O11").Value = 145.3

$ws.Range("O12").PasteSpecial(-4122)
$ws.Range("O12").Value = 98.8

$ws.Range("O13").PasteSpecial(-4122)
$ws.Range("O13").Value = 98.7

$ws.Range("O14").PasteSpecial(-4122)
$ws.Range("O14").Value = 1.8

# --- update selection to match the authored change ---
$ws.Range("P5").Select()
